$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114
$ws.Range("B114").Value = 6838437
$ws.Range("F114").Value = "FC Twente"
$ws.Range("G114").Value = "PSV"
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 3
$ws.Range("J114").Value = "A"
$ws.Range("K114").Value = 4
$ws.Range("M114").Value = 1.833
$ws.Range("N114").Value = 3.6
$ws.Range("O114").Value = 3.6
$ws.Range("P114").Value = 2
$ws.Range("Q114").Value = 0.5
$ws.Range("R114").Value = 1.825
$ws.Range("S114").Value = 2.025
$ws.Range("T114").Value = 3
$ws.Range("U114").Value = 2.025
$ws.Range("V114").Value = 1.825
$ws.Range("W114").Value = -1
$ws.Range("Y114").Value = 1
$ws.Range("Z114").Value = -1
$ws.Range("AA114").Value = 1.025
$ws.Range("AB114").Value = 0
$ws.Range("AC114").Value = -0

# Row 115
$ws.Range("B115").Value = 6838438
$ws.Range("F115").Value = "Heerenveen"
$ws.Range("G115").Value = "Fortuna Sittard"
$ws.Range("H115").Value = 3
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = "H"
$ws.Range("K115").Value = 1.833
$ws.Range("M115").Value = 4
$ws.Range("N115").Value = 1.95
$ws.Range("O115").Value = 3.5
$ws.Range("P115").Value = 3.8
$ws.Range("Q115").Value = -0.5
$ws.Range("R115").Value = 1.975
$ws.Range("S115").Value = 1.875
$ws.Range("T115").Value = 2.5
$ws.Range("U115").Value = 2
$ws.Range("V115").Value = 1.85
$ws.Range("W115").Value = 0.95
$ws.Range("Y115").Value = -1
$ws.Range("Z115").Value = 0.9750000000000001
$ws.Range("AA115").Value = -1
$ws.Range("AB115").Value = 1
$ws.Range("AC115").Value = -1

# Row 116
$ws.Range("B116").Value = 6838440
$ws.Range("F116").Value = "Ajax"
$ws.Range("G116").Value = "Vitesse"
$ws.Range("H116").Value = 5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = "H"
$ws.Range("K116").Value = 1.363
$ws.Range("L116").Value = 5
$ws.Range("M116").Value = 7.5
$ws.Range("N116").Value = 1.333
$ws.Range("O116").Value = 6
$ws.Range("P116").Value = 7.5
$ws.Range("Q116").Value = -1.5
$ws.Range("R116").Value = 1.825
$ws.Range("S116").Value = 2.025
$ws.Range("T116").Value = 3.5
$ws.Range("U116").Value = 1.95
$ws.Range("V116").Value = 1.9
$ws.Range("W116").Value = 0.333
$ws.Range("Y116").Value = -1
$ws.Range("Z116").Value = 0.825
$ws.Range("AA116").Value = -1
$ws.Range("AB116").Value = 0.95
$ws.Range("AC116").Value = -1

# Row 117
$ws.Range("B117").Value = 6838439
$ws.Range("F117").Value = "PEC Zwolle"
$ws.Range("G117").Value = "RKC"
$ws.Range("H117").Value = 1
$ws.Range("I117").Value = 2
$ws.Range("J117").Value = "A"
$ws.Range("K117").Value = 1.909
$ws.Range("L117").Value = 3.6
$ws.Range("M117").Value = 3.75
$ws.Range("N117").Value = 2.05
$ws.Range("O117").Value = 3.6
$ws.Range("P117").Value = 3.4
$ws.Range("Q117").Value = -0.5
$ws.Range("R117").Value = 2.05
$ws.Range("S117").Value = 1.8
$ws.Range("T117").Value = 2.75
$ws.Range("U117").Value = 1.9
$ws.Range("V117").Value = 1.95
$ws.Range("W117").Value = -1
$ws.Range("Y117").Value = 2.4
$ws.Range("Z117").Value = -1
$ws.Range("AA117").Value = 0.8
$ws.Range("AB117").Value = 0.45
$ws.Range("AC117").Value = -0.5

# Row 251
$ws.Range("B251").Value = 6973370
$ws.Range("E251").Value = 45388.47916666666
$ws.Range("F251").Value = "Sparta Rotterdam"
$ws.Range("G251").Value = "Heracles"
$ws.Range("K251").Value = 1.75
$ws.Range("L251").Value = 3.8
$ws.Range("M251").Value = 4.2
$ws.Range("N251").Value = 1.666
$ws.Range("O251").Value = 4
$ws.Range("P251").Value = 4.5
$ws.Range("Q251").Value = -0.75
$ws.Range("R251").Value = 1.89
$ws.Range("S251").Value = 2.01
$ws.Range("T251").Value = 3
$ws.Range("U251").Value = 2.05
$ws.Range("V251").Value = 1.8

# Row 252
$ws.Range("B252").Value = 6838570
$ws.Range("E252").Value = 45388.57291666666
$ws.Range("F252").Value = "PSV"
$ws.Range("G252").Value = "AZ"
$ws.Range("K252").Value = 1.4
$ws.Range("L252").Value = 4.75
$ws.Range("M252").Value = 7
$ws.Range("N252").Value = 1.45
$ws.Range("O252").Value = 4.5
$ws.Range("P252").Value = 6
$ws.Range("Q252").Value = -1.25
$ws.Range("R252").Value = 2.01
$ws.Range("S252").Value = 1.89
$ws.Range("T252").Value = 3
$ws.Range("U252").Value = 2
$ws.Range("V252").Value = 1.85

# Row 253
$ws.Range("B253").Value = 6956565
$ws.Range("E253").Value = 45388.625
$ws.Range("F253").Value = "PEC Zwolle"
$ws.Range("G253").Value = "Excelsior"
$ws.Range("K253").Value = 2.15
$ws.Range("L253").Value = 3.6
$ws.Range("M253").Value = 3.1
$ws.Range("N253").Value = 2
$ws.Range("O253").Value = 3.6
$ws.Range("P253").Value = 3.4
$ws.Range("Q253").Value = -0.5
$ws.Range("R253").Value = 2.04
$ws.Range("S253").Value = 1.86
$ws.Range("T253").Value = 3

# Row 254
$ws.Range("B254").Value = 6956849
$ws.Range("E254").Value = 45388.66666666666
$ws.Range("F254").Value = "FC Twente"
$ws.Range("G254").Value = "Fortuna Sittard"
$ws.Range("K254").Value = 1.3
$ws.Range("L254").Value = 5.25
$ws.Range("M254").Value = 9
$ws.Range("N254").Value = 1.4
$ws.Range("O254").Value = 4.75
$ws.Range("P254").Value = 6.5
$ws.Range("Q254").Value = -1.25
$ws.Range("R254").Value = 1.99
$ws.Range("S254").Value = 1.91
$ws.Range("U254").Value = 1.975
$ws.Range("V254").Value = 1.875

# Row 255
$ws.Range("B255").Value = 6838586
$ws.Range("E255").Value = 45389.30208333334
$ws.Range("F255").Value = "Vitesse"
$ws.Range("G255").Value = "NEC"
$ws.Range("K255").Value = 2.8
$ws.Range("L255").Value = 3.6
$ws.Range("M255").Value = 2.3
$ws.Range("N255").Value = 2.875
$ws.Range("O255").Value = 3.6
$ws.Range("P255").Value = 2.25
$ws.Range("Q255").Value = 0.25
$ws.Range("R255").Value = 1.9
$ws.Range("S255").Value = 2
$ws.Range("T255").Value = 2.75
$ws.Range("U255").Value = 1.925
$ws.Range("V255").Value = 1.925

# Row 256
$ws.Range("B256").Value = 6994877
$ws.Range("E256").Value = 45389.39583333334
$ws.Range("F256").Value = "Feyenoord"
$ws.Range("G256").Value = "Ajax"
$ws.Range("K256").Value = 1.6
$ws.Range("L256").Value = 4
$ws.Range("M256").Value = 5.25
$ws.Range("N256").Value = 1.45
$ws.Range("O256").Value = 4.5
$ws.Range("P256").Value = 6.5
$ws.Range("Q256").Value = -1.25
$ws.Range("R256").Value = 2
$ws.Range("S256").Value = 1.9
$ws.Range("T256").Value = 3.25
$ws.Range("U256").Value = 2.025
$ws.Range("V256").Value = 1.825

# Row 257
$ws.Range("B257").Value = 6838569
$ws.Range("E257").Value = 45389.39583333334
$ws.Range("F257").Value = "Go Ahead Eagles"
$ws.Range("G257").Value = "Almere City FC"
$ws.Range("K257").Value = 1.75
$ws.Range("L257").Value = 3.6
$ws.Range("M257").Value = 4.5
$ws.Range("N257").Value = 1.75
$ws.Range("O257").Value = 3.6
$ws.Range("P257").Value = 4.333
$ws.Range("Q257").Value = -0.75
$ws.Range("R257").Value = 2.06
$ws.Range("S257").Value = 1.84
$ws.Range("U257").Value = 2.025
$ws.Range("V257").Value = 1.825

# Row 258
$ws.Range("B258").Value = 6920184
$ws.Range("E258").Value = 45389.48958333334
$ws.Range("F258").Value = "FC Volendam"
$ws.Range("G258").Value = "RKC"
$ws.Range("K258").Value = 2.7
$ws.Range("L258").Value = 3.5
$ws.Range("M258").Value = 2.4
$ws.Range("N258").Value = 3.75
$ws.Range("O258").Value = 3.75
$ws.Range("P258").Value = 1.909
$ws.Range("Q258").Value = 0.5
$ws.Range("R258").Value = 1.97
$ws.Range("S258").Value = 1.93
$ws.Range("U258").Value = 1.875
$ws.Range("V258").Value = 1.975

# Row 259
$ws.Range("N259").Value = 2.55
$ws.Range("P259").Value = 2.5
$ws.Range("R259").Value = 1.99
$ws.Range("S259").Value = 1.91

